$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'61.859.14"
$ws.Range('E2').Value = '  +1.47%  '
$ws.Range('D3').Value = "'3.416.07"
$ws.Range('E3').Value = '  +4.19%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = "'577.28"
$ws.Range('E5').Value = '  +2.64%  '
$ws.Range('D6').Value = "'139.26"
$ws.Range('E6').Value = '  +10.26%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = "'3.415.27"
$ws.Range('E8').Value = '  +4.20%  '
$ws.Range('E9').Value = '  +0.74%  '
$ws.Range('D10').Value = "'7.70"
$ws.Range('E10').Value = '  +6.17%  '
$ws.Range('D11').Value = "'0.128"
$ws.Range('E11').Value = '  +9.25%  '
$ws.Range('E12').Value = '  +6.41%  '
$ws.Range('D13').Value = "'3.994.26"
$ws.Range('E13').Value = '  +3.63%  '
$ws.Range('E14').Value = '  +2.16%  '
$ws.Range('E15').Value = '  +9.29%  '
$ws.Range('D16').Value = "'3.416.32"
$ws.Range('E16').Value = '  +3.88%  '
$ws.Range('D17').Value = "'25.65"
$ws.Range('E17').Value = '  +7.32%  '
$ws.Range('D18').Value = "'61.885.14"
$ws.Range('E18').Value = '  +1.24%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = "'14.13"
$ws.Range('E19').Value = '  +7.28%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').Value = "'5.96"
$ws.Range('E20').Value = '  +7.30%  '
$ws.Range('D21').Value = "'9.48"
$ws.Range('E21').Value = '  +7.09%  '
$ws.Range('D22').Value = "'391.59"
$ws.Range('E22').Value = '  +11.97%  '
$ws.Range('D23').Value = "'0.574"
$ws.Range('E23').Value = '  +4.82%  '
$ws.Range('D24').Value = "'3.551.79"
$ws.Range('E24').Value = '  +3.94%  '
$ws.Range('E25').Value = '  +0.15%  '
$ws.Range('E26').Value = '  +20.60%  '
$ws.Range('D27').Value = "'71.35"
$ws.Range('E27').Value = '  +3.67%  '
$ws.Range('E28').Value = '  +16.61%  '
$ws.Range('D29').Value = "'7.81"
$ws.Range('E29').Value = '  +11.08%  '
$ws.Range('D30').Value = "'0.993"
$ws.Range('E30').Value = '  -0.59%  '
$ws.Range('D31').Value = "'8.33"
$ws.Range('E31').Value = '  +8.13%  '
$ws.Range('E32').Value = '  +8.27%  '
$ws.Range('E33').Value = '  +3.63%  '
$ws.Range('D34').Value = "'3.448.97"
$ws.Range('E34').Value = '  +4.17%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').Value = "'23.66"
$ws.Range('E36').Value = '  +5.64%  '
$ws.Range('D37').Value = "'5.54"
$ws.Range('E37').Value = '  +6.95%  '
$ws.Range('D38').Value = "'7.06"
$ws.Range('E38').Value = '  +5.15%  '
$ws.Range('E39').Value = '  +7.19%  '
$ws.Range('D40').Value = "'161.92"
$ws.Range('E40').Value = '  -0.78%  '
$ws.Range('D41').Value = "'0.0801"
$ws.Range('E41').Value = '  +7.41%  '
$ws.Range('E42').Value = '  +12.99%  '
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('E44').Value = '  +4.51%  '
$ws.Range('E45').Value = '  +10.74%  '
$ws.Range('D46').Value = "'0.777"
$ws.Range('E46').Value = '  +5.83%  '
$ws.Range('D47').Value = "'41.25"
$ws.Range('E47').Value = '  +0.66%  '
$ws.Range('D48').Value = "'23.50"
$ws.Range('E48').Value = '  +7.14%  '
$ws.Range('D49').Value = "'7.02"
$ws.Range('E49').Value = '  +6.11%  '
$ws.Range('D50').Value = "'22.91"
$ws.Range('E50').Value = '  +9.49%  '
$ws.Range('D51').Value = "'2.364.28"
$ws.Range('E51').Value = '  +10.62%  '

# Strip the quote-prefix formatting introduced by forcing text values,
# restoring cells to the default (unstyled) appearance.
$ws.Range('D2').ClearFormats()
$ws.Range('D3').ClearFormats()
$ws.Range('D5').ClearFormats()
$ws.Range('D6').ClearFormats()
$ws.Range('D8').ClearFormats()
$ws.Range('D10').ClearFormats()
$ws.Range('D11').ClearFormats()
$ws.Range('D13').ClearFormats()
$ws.Range('D16').ClearFormats()
$ws.Range('D17').ClearFormats()
$ws.Range('D18').ClearFormats()
$ws.Range('D19').ClearFormats()
$ws.Range('D20').ClearFormats()
$ws.Range('D21').ClearFormats()
$ws.Range('D22').ClearFormats()
$ws.Range('D23').ClearFormats()
$ws.Range('D24').ClearFormats()
$ws.Range('D27').ClearFormats()
$ws.Range('D29').ClearFormats()
$ws.Range('D30').ClearFormats()
$ws.Range('D31').ClearFormats()
$ws.Range('D34').ClearFormats()
$ws.Range('D36').ClearFormats()
$ws.Range('D37').ClearFormats()
$ws.Range('D38').ClearFormats()
$ws.Range('D40').ClearFormats()
$ws.Range('D41').ClearFormats()
$ws.Range('D46').ClearFormats()
$ws.Range('D47').ClearFormats()
$ws.Range('D48').ClearFormats()
$ws.Range('D49').ClearFormats()
$ws.Range('D50').ClearFormats()
$ws.Range('D51').ClearFormats()
